# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (date 2021-11-09 / serial 44509) for
# "Alcachofa Española" (Extra/Primera/Segunda) at the top of the existing
# date-ordered block, pushing the rest of the sheet down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 467 (shifts
# everything from 467 downwards to 470 onwards).
$ws.Rows("467:469").Insert()

# --- Row 467: Española / Extra ---
$ws.Range("A467").Value = 6
$ws.Range("B467").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C467").Value = "Metropolitana"
$ws.Range("D467").Value = 44509
$ws.Range("E467").Value = 13
$ws.Range("F467").Value = 100112013
$ws.Range("G467").Value = "Alcachofa"
$ws.Range("H467").Value = "Española"
$ws.Range("I467").Value = "Extra"
$ws.Range("J467").Value = 3800
$ws.Range("K467").Value = 300
$ws.Range("L467").Value = 320
$ws.Range("M467").Value = 307
$ws.Range("N467").Value = "$/unidad"
$ws.Range("O467").Value = "Región Metropolitana"
$ws.Range("P467").Value = 307
$ws.Range("Q467").Value = 1
$ws.Range("R467").Value = "Hortaliza"

# --- Row 468: Española / Primera ---
$ws.Range("A468").Value = 6
$ws.Range("B468").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C468").Value = "Metropolitana"
$ws.Range("D468").Value = 44509
$ws.Range("E468").Value = 13
$ws.Range("F468").Value = 100112013
$ws.Range("G468").Value = "Alcachofa"
$ws.Range("H468").Value = "Española"
$ws.Range("I468").Value = "Primera"
$ws.Range("J468").Value = 3000
$ws.Range("K468").Value = 250
$ws.Range("L468").Value = 270
$ws.Range("M468").Value = 260
$ws.Range("N468").Value = "$/unidad"
$ws.Range("O468").Value = "Región Metropolitana"
$ws.Range("P468").Value = 260
$ws.Range("Q468").Value = 1
$ws.Range("R468").Value = "Hortaliza"

# --- Row 469: Española / Segunda ---
$ws.Range("A469").Value = 6
$ws.Range("B469").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C469").Value = "Metropolitana"
$ws.Range("D469").Value = 44509
$ws.Range("E469").Value = 13
$ws.Range("F469").Value = 100112013
$ws.Range("G469").Value = "Alcachofa"
$ws.Range("H469").Value = "Española"
$ws.Range("I469").Value = "Segunda"
$ws.Range("J469").Value = 1900
$ws.Range("K469").Value = 200
$ws.Range("L469").Value = 230
$ws.Range("M469").Value = 214
$ws.Range("N469").Value = "$/unidad"
$ws.Range("O469").Value = "Región Metropolitana"
$ws.Range("P469").Value = 214
$ws.Range("Q469").Value = 1
$ws.Range("R469").Value = "Hortaliza"
